$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M0 - Account Mgmt")

# Row 17: add "Done" status cell (C17) - reuses existing "Good" cell style
$ws.Range("C17").Value = "Done"
$ws.Range("C17").Style = "Good"

# Set new text values in the same order the author typed them, so the
# shared-string table grows in the same order as the target workbook.
$ws.Range("C18").Value = "NEXT"
$ws.Range("D20").Value = "Requires some design. End goal is to support multiple games. Easier to do now rather than later.  Should give it a better name than UserGames. E.g. GameLogs"
$ws.Range("C20").Value = "REVISIT"
$ws.Range("D18").Value = "Ties to item below. "

# Row 19 reuses the "NEXT" marker created above
$ws.Range("C19").Value = "NEXT"

# Apply bold formatting
$ws.Range("C18").Font.Bold = $true
$ws.Range("C19").Font.Bold = $true
$ws.Range("A20").Font.Bold = $true
$ws.Range("B20").Font.Bold = $true
$ws.Range("C20").Font.Bold = $true
$ws.Range("D20").Font.Bold = $true
$ws.Range("E20").Font.Bold = $true

# Update selection to match where the edits were made
$ws.Activate()
$ws.Range("C20").Select()
